$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Платежи")

$ws.Range("B2").Value = "2011-Sep-13 / 15:09"
$ws.Range("B3").Value = "2011-Sep-13 / 15:09"
